$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header text updates (rich-text shared strings)
# ---------------------------------------------------------------------------
# "Volume 32   Number  8" -> "...Number  9"
$ws.Range("A8").Characters(21, 1).Text = "9"

# "Report Covering the Week  2/17/2025  Through  2/23/2025"
#   -> "...2/24/2025  Through  3/2/2025"
$ws.Range("C9").Characters(27, 9).Text = "2/24/2025"
$ws.Range("C9").Characters(47, 9).Text = "3/2/2025"

# ---------------------------------------------------------------------------
# 2. Cells that flip between text placeholder ("0" / "***.*") and a real
#    number need both their value AND their number-format style changed.
#    Trick: assign the value first, then copy the *format only* from a
#    neighbouring cell that already carries the desired numeric style.
# ---------------------------------------------------------------------------

# Row 15
$ws.Range("C15").Value = 1
$ws.Range("D15").Copy()
$ws.Range("C15").PasteSpecial(-4122)

# Row 22
$ws.Range("D22").Value = 1
$ws.Range("F22").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").Value = -100
$ws.Range("H22").Copy()
$ws.Range("E22").PasteSpecial(-4122)

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("D27").Copy()
$ws.Range("C27").PasteSpecial(-4122)

# Row 28
$ws.Range("D28").Value = 1
$ws.Range("C28").Copy()
$ws.Range("D28").PasteSpecial(-4122)

$ws.Range("E28").Value = 100
$ws.Range("H28").Copy()
$ws.Range("E28").PasteSpecial(-4122)

# Row 29-31: F column flips the other way, number -> text "0"
foreach ($addr in @("F29", "F30", "F31")) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = "0"
}
$ws.Range("C22").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("F31").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Plain numeric value updates (style unchanged)
# ---------------------------------------------------------------------------

# Row 15
$ws.Range("E15").Value = 0
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -33.333333333333
$ws.Range("I15").Value = 4
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 33.333333333333
$ws.Range("L15").Value = 33.333333333333
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = -42.857142857142

# Row 16
$ws.Range("D16").Value = 3
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -83.333333333333
$ws.Range("J16").Value = 36
$ws.Range("K16").Value = -61.111111111111
$ws.Range("L16").Value = -61.111111111111
$ws.Range("M16").Value = -63.157894736842
$ws.Range("N16").Value = -92.964824120603

# Row 17
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 23.809523809523
$ws.Range("I17").Value = 46
$ws.Range("J17").Value = 47
$ws.Range("K17").Value = -2.127659574468
$ws.Range("L17").Value = -11.538461538461
$ws.Range("M17").Value = 53.333333333333
$ws.Range("N17").Value = 43.75

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 11
$ws.Range("E18").Value = -72.727272727272
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -37.5
$ws.Range("I18").Value = 25
$ws.Range("J18").Value = 37
$ws.Range("K18").Value = -32.432432432432
$ws.Range("L18").Value = -32.432432432432
$ws.Range("M18").Value = -68.75
$ws.Range("N18").Value = -92.647058823529

# Row 19
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = -29.824561403508
$ws.Range("I19").Value = 79
$ws.Range("J19").Value = 137
$ws.Range("K19").Value = -42.335766423357
$ws.Range("L19").Value = -24.761904761904
$ws.Range("M19").Value = 16.176470588235
$ws.Range("N19").Value = -21

# Row 20
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = -62.5
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = -79.310344827586
$ws.Range("I20").Value = 26
$ws.Range("J20").Value = 61
$ws.Range("K20").Value = -57.377049180327
$ws.Range("L20").Value = -50
$ws.Range("M20").Value = -29.729729729729
$ws.Range("N20").Value = -95.786061588330

# Row 21 (TOTAL)
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 45
$ws.Range("F21").Value = 92
$ws.Range("G21").Value = 152
$ws.Range("H21").Value = -39.473684210526
$ws.Range("I21").Value = 194
$ws.Range("J21").Value = 321
$ws.Range("K21").Value = -39.563862928348
$ws.Range("L21").Value = -32.167832167832
$ws.Range("M21").Value = -23.921568627451
$ws.Range("N21").Value = -85.053929121725

# Row 22 (remaining plain numeric updates)
$ws.Range("J22").Value = 5
$ws.Range("K22").Value = -60
$ws.Range("M22").Value = -60

# Row 24
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -31.034482758620
$ws.Range("F24").Value = 80
$ws.Range("G24").Value = 130
$ws.Range("H24").Value = -38.461538461538
$ws.Range("I24").Value = 200
$ws.Range("J24").Value = 267
$ws.Range("K24").Value = -25.093632958801
$ws.Range("L24").Value = -9.502262443438
$ws.Range("M24").Value = 11.731843575419

# Row 25
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 11.111111111111
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = 8.823529411764
$ws.Range("I25").Value = 83
$ws.Range("J25").Value = 76
$ws.Range("K25").Value = 9.210526315789
$ws.Range("L25").Value = 40.677966101694

# Row 26
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 32
$ws.Range("G26").Value = 51
$ws.Range("H26").Value = -37.254901960784
$ws.Range("I26").Value = 77
$ws.Range("J26").Value = 90
$ws.Range("K26").Value = -14.444444444444
$ws.Range("L26").Value = -6.097560975609
$ws.Range("M26").Value = -36.363636363636

# Row 27 (remaining plain numeric updates)
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 5
$ws.Range("J27").Value = 4
$ws.Range("K27").Value = 25
$ws.Range("L27").Value = 25

# Row 28 (remaining plain numeric updates)
$ws.Range("C28").Value = 2
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 250
$ws.Range("I28").Value = 9
$ws.Range("J28").Value = 4
$ws.Range("K28").Value = 125
$ws.Range("L28").Value = -18.181818181818

# Row 29-31 (remaining plain numeric updates)
$ws.Range("H29").Value = -100
$ws.Range("H30").Value = -100
$ws.Range("H31").Value = -100
